$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-27 Saturday" "2025-09-28 Sunday"

Replace-Text "272×9=" "926×3="
Replace-Text "199×2=" "175×5="
Replace-Text "165×3=" "541×4="
Replace-Text "244×2=" "766×8="
Replace-Text "732×3=" "948×8="

Replace-Text "636×8=" "554×4="
Replace-Text "356×9=" "897×7="
Replace-Text "566×7=" "812×4="
Replace-Text "326×9=" "118×6="
Replace-Text "149×3=" "767×3="

Replace-Text "332×4=" "852×5="
Replace-Text "768×5=" "904×8="
Replace-Text "493×7=" "247×7="
Replace-Text "114×9=" "203×8="
Replace-Text "804×4=" "742×4="

Replace-Text "293×8=" "995×9="
Replace-Text "547×3=" "707×3="
Replace-Text "194×7=" "405×5="
Replace-Text "329×2=" "146×6="
Replace-Text "615×2=" "417×5="

Replace-Text "430×7=" "760×7="
Replace-Text "808×5=" "650×4="
Replace-Text "852×3=" "864×3="
Replace-Text "979×2=" "994×4="
Replace-Text "671×3=" "693×7="
